$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''43.895.46'
$ws.Range('E2').Value = '  -0.17%  '

$ws.Range('D3').Value = '''2.229.23'
$ws.Range('E3').Value = '  -1.51%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.25%  '

$ws.Range('D5').Value = '''301.92'
$ws.Range('E5').Value = '  -5.23%  '

$ws.Range('D6').Value = '''93.59'
$ws.Range('E6').Value = '  -7.91%  '

$ws.Range('D7').Value = '''0.565'
$ws.Range('E7').Value = '  -1.93%  '

$ws.Range('E8').Value = '  +0.15%  '

$ws.Range('D9').Value = '''0.515'
$ws.Range('E9').Value = '  -7.29%  '

$ws.Range('D10').Value = '''34.02'
$ws.Range('E10').Value = '  -8.51%  '

$ws.Range('D11').Value = '''0.0796'
$ws.Range('E11').Value = '  -4.40%  '

$ws.Range('D12').Value = '''7.05'
$ws.Range('E12').Value = '  -7.81%  '

$ws.Range('D13').Value = '''0.103'
$ws.Range('E13').Value = '  -3.24%  '

$ws.Range('D14').Value = '''2.566.85'
$ws.Range('E14').Value = '  -1.59%  '

$ws.Range('D15').Value = '''2.257.85'
$ws.Range('E15').Value = '  -0.42%  '

$ws.Range('D16').Value = '''0.806'
$ws.Range('E16').Value = '  -6.65%  '

$ws.Range('D17').Value = '''13.35'
$ws.Range('E17').Value = '  -7.78%  '

$ws.Range('D18').Value = '''43.655.35'
$ws.Range('E18').Value = '  -0.51%  '

$ws.Range('D19').Value = '''0.0₃0944'
$ws.Range('E19').Value = '  -4.20%  '

$ws.Range('D20').Value = '''11.94'
$ws.Range('E20').Value = '  -10.69%  '

$ws.Range('D21').Value = '''6.09'
$ws.Range('E21').Value = '  -7.28%  '

$ws.Range('D22').Value = '''64.07'
$ws.Range('E22').Value = '  -2.47%  '

$ws.Range('D23').Value = '''234.35'
$ws.Range('E23').Value = '  -0.39%  '

$ws.Range('D24').Value = '''2.88'
$ws.Range('E24').Value = '  -8.05%  '

$ws.Range('E25').Value = '  -0.03%  '

$ws.Range('D26').Value = '''1.91'
$ws.Range('E26').Value = '  -9.66%  '

$ws.Range('D27').Value = '''9.72'
$ws.Range('E27').Value = '  -4.18%  '

$ws.Range('E28').Value = '  -2.62%  '

$ws.Range('D29').Value = '''35.91'
$ws.Range('E29').Value = '  -3.52%  '

$ws.Range('D30').Value = '''5.83'
$ws.Range('E30').Value = '  -6.51%  '

$ws.Range('D31').Value = '''19.69'
$ws.Range('E31').Value = '  -2.72%  '

$ws.Range('D32').Value = '''151.35'
$ws.Range('E32').Value = '  -4.39%  '

$ws.Range('D33').Value = '''0.0796'
$ws.Range('E33').Value = '  -6.86%  '

$ws.Range('D34').Value = '''3.21'
$ws.Range('E34').Value = '  +4.40%  '

$ws.Range('D35').Value = '''2.60'
$ws.Range('E35').Value = '  -4.36%  '

$ws.Range('D36').Value = '''0.117'
$ws.Range('E36').Value = '  -1.99%  '

$ws.Range('E37').Value = '  -8.41%  '

$ws.Range('D38').Value = '''1.74'
$ws.Range('E38').Value = '  -11.12%  '

$ws.Range('D39').Value = '''14.43'
$ws.Range('E39').Value = '  -10.98%  '

$ws.Range('D40').Value = '''3.77'
$ws.Range('E40').Value = '  -10.50%  '

$ws.Range('D41').Value = '''3.24'
$ws.Range('E41').Value = '  -12.90%  '

$ws.Range('D42').Value = '''0.0294'
$ws.Range('E42').Value = '  -6.97%  '

$ws.Range('E43').Value = '  +0.26%  '

$ws.Range('D44').Value = '''1.727.25'
$ws.Range('E44').Value = '  -4.08%  '

$ws.Range('D45').Value = '''82.65'
$ws.Range('E45').Value = '  +0.24%  '

$ws.Range('D46').Value = '''4.89'
$ws.Range('E46').Value = '  -6.31%  '

$ws.Range('D47').Value = '''0.183'
$ws.Range('E47').Value = '  -7.77%  '

$ws.Range('D48').Value = '''98.43'
$ws.Range('E48').Value = '  -6.27%  '

$ws.Range('D49').Value = '''7.97'
$ws.Range('E49').Value = '  -4.81%  '

$ws.Range('D50').Value = '''67.47'
$ws.Range('E50').Value = '  -11.52%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''12.54'
$ws.Range('E51').Value = '  -9.47%  '
